# Translate the "ContosoLearn Competitor SWOT" document from English to
# Spanish. Uses Find/Replace (ReplaceAll) over the whole document content
# so that formatting on the surrounding runs (bold SWOT labels, bullet
# numbering, etc.) is preserved while only the visible text changes.

$d = $word.ActiveDocument

function Replace-All($findText, $replaceText) {
    $d.Content.Find.Execute($findText, $false, $false, $false, $false, $false, `
                             $true, 1, $false, $replaceText, 2) | Out-Null
}

# Title
Replace-All "ContosoLearn Competitor SWOT" "DAFO de competidores ContosoLearn"

# Bold SWOT labels (shared across both companies' sections)
Replace-All "Strengths:" "Fortalezas:"
Replace-All "Weaknesses:" "Puntos débiles:"
Replace-All "Opportunities:" "Oportunidades:"
Replace-All "Threats:" "Amenazas:"

# Fabrikam Learning - Strengths
Replace-All " Fabrikam Learning provides a comprehensive set of analytics and reporting tools. It ensures the continuous monitoring of teaching and learning activities, as well as pinpointing problematic areas that need to be addressed." `
            " Fabrikam Learning proporciona un conjunto completo de herramientas de análisis e informes. Garantiza la supervisión continua de las actividades de enseñanza y aprendizaje, así como la identificación de áreas problemáticas que deben abordarse."

# Fabrikam Learning - Weaknesses
Replace-All " While Fabrikam Learning has robust reporting capabilities, it might be overwhelming for some users due to its comprehensive nature." `
            " aunque Fabrikam Learning tiene funcionalidades de informes sólidas, puede ser abrumador para algunos usuarios debido a su naturaleza completa."

# Fabrikam Learning - Opportunities
Replace-All " There is a growing demand for personalized learning experiences and data-driven recommendations. Fabrikam Learning can leverage its robust analytics and reporting tools to meet this demand." `
            " hay una creciente demanda de experiencias de aprendizaje personalizadas y recomendaciones controladas por datos. Fabrikam Learning puede aprovechar sus sólidas herramientas de análisis e informes para satisfacer esta demanda."

# Fabrikam Learning - Threats
Replace-All " The eLearning market is highly competitive with many players offering similar features. Fabrikam Learning needs to continuously innovate to stay ahead." `
            " el mercado de eLearning es altamente competitivo con muchos jugadores que ofrecen características similares. Fabrikam Learning debe innovar continuamente para mantenerse a la vanguardia."

# AdatumLearn - Strengths
Replace-All " AdatumLearn offers courses on business analysis techniques such as MOST and SWOT. This shows their commitment to providing valuable content to their users." `
            " AdatumLearn ofrece cursos sobre técnicas de análisis de negocios como MOST y DAFO. Esto muestra su compromiso de proporcionar contenido valioso a sus usuarios."

# AdatumLearn - Weaknesses
Replace-All " The information provided in their courses is a compilation of third-party generated information. This might not be as valuable as original content." `
            " la información proporcionada en sus cursos es una compilación de información generada por terceros. Esto podría no ser tan valioso como el contenido original."

# AdatumLearn - Opportunities
Replace-All " AdatumLearn can create more original content to provide unique value to their users. They can also expand their course offerings to cover more topics." `
            " AdatumLearn puede crear contenido más original para proporcionar un valor único a sus usuarios. También puede ampliar sus ofertas de cursos para tratar más temas."

# AdatumLearn - Threats (ends with a straight double quote, matching the
# source). Replace the bulk of the text first...
Replace-All 'Like Fabrikam Learning, AdatumLearn also faces stiff competition in the eLearning market. They need to continuously improve their offerings to stay competitive."' `
            'al igual que Fabrikam Learning, AdatumLearn también se enfrenta a una competencia rígida en el mercado de eLearning. Necesita mejorar continuamente su oferta para mantenerse competitivo".'

# ...then repair the trailing quote: Find/Replace's text-insertion runs
# Word's AutoCorrect "smart quotes", which turns the straight `"` above
# into a curly right double quotation mark. Locate that curly quote and
# swap it back for a literal straight quote via direct range surgery
# (Delete + InsertAfter bypasses AutoCorrect, unlike Find.Execute).
$curlyQuote = [char]0x201D
$straightQuote = [char]0x22
$hit = $d.Content
$found = $hit.Find.Execute($curlyQuote, $false, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)
if ($found) {
    $hit.Delete()
    $hit.InsertAfter($straightQuote)
}
